$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column cells hold plain text values (e.g. "68.576.23",
# "1.00", "8.27"); some look like plain numbers, so force the Text
# number format on those specific cells before writing so Excel does
# not silently coerce them into numeric values (which would drop
# formatting such as trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.548.83'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.697.91'
$ws.Range("E3").Value = '  +2.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.67'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.93'
$ws.Range("E6").Value = '  +2.02%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.698.51'
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("E10").Value = '  -2.79%  '
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("E13").Value = '  +2.40%  '
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.194.10'
$ws.Range("E15").Value = '  +2.22%  '
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.500.99'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.696.84'
$ws.Range("E18").Value = '  +2.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.96'
$ws.Range("E19").Value = '  +5.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.68'
$ws.Range("E20").Value = '  +3.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '366.55'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.55'
$ws.Range("E22").Value = '  +2.88%  '
$ws.Range("E23").Value = '  +1.89%  '
$ws.Range("E24").Value = '  +2.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.47'
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  +4.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.835.47'
$ws.Range("E28").Value = '  +2.05%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '572.83'
$ws.Range("E31").Value = '  +3.18%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.46'
$ws.Range("E32").Value = '  +3.98%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.23'
$ws.Range("E33").Value = '  +2.66%  '
$ws.Range("E34").Value = '  +5.50%  '
$ws.Range("E35").Value = '  +3.01%  '
$ws.Range("E36").Value = '  +5.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.97'
$ws.Range("E38").Value = '  +3.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '161.01'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.380'
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("E41").Value = '  +1.69%  '
$ws.Range("E42").Value = '  +1.78%  '
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0317'
$ws.Range("E46").Value = '  -6.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.42'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("E48").Value = '  +5.89%  '
$ws.Range("E49").Value = '  +5.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.602'
$ws.Range("E50").Value = '  +7.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.07'
$ws.Range("E51").Value = '  +0.04%  '
